$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.322531
$ws.Range("H2").Value = 0.967593
$ws.Range("I2").Value = 0.01892149513432853
$ws.Range("J2").Value = 0.01892149513432853
$ws.Range("M2").Value = 0.004862
$ws.Range("P2").Value = 0.02348134339170667
$ws.Range("Q2").Value = 0.001568145722
$ws.Range("R2").Value = 0.014113311498
$ws.Range("S2").Value = 0.0004443021247336752
$ws.Range("T2").Value = 0.0004443021247336752
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.322531
$ws.Range("H3").Value = 0.967593
$ws.Range("I3").Value = 0.01892149513432853
$ws.Range("J3").Value = 0.01892149513432853
$ws.Range("Q3").Value = 0.065214478076
$ws.Range("R3").Value = 0.586930302684
$ws.Range("S3").Value = 0.01847719300959486
$ws.Range("T3").Value = 0.01847719300959486
$ws.Range("I4").Value = 0.0261208867009986
$ws.Range("J4").Value = 0.0261208867009986
$ws.Range("M4").Value = 0.004862
$ws.Range("P4").Value = 0.02348134339170667
$ws.Range("S4").Value = 0.000613353510322012
$ws.Range("T4").Value = 0.000613353510322012
$ws.Range("I5").Value = 0.0261208867009986
$ws.Range("J5").Value = 0.0261208867009986
$ws.Range("S5").Value = 0.02550753319067658
$ws.Range("T5").Value = 0.02550753319067658
$ws.Range("I6").Value = 0.954957618164673
$ws.Range("J6").Value = 0.954957618164673
$ws.Range("M6").Value = 0.004862
$ws.Range("P6").Value = 0.02348134339170667
$ws.Range("Q6").Value = 0.07914346583
$ws.Range("R6").Value = 0.71229119247
$ws.Range("S6").Value = 0.02242368775665098
$ws.Range("T6").Value = 0.02242368775665099
$ws.Range("I7").Value = 0.954957618164673
$ws.Range("J7").Value = 0.954957618164673
$ws.Range("S7").Value = 0.9325339304080219
$ws.Range("T7").Value = 0.932533930408022
